$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is inserted at row 100, pushing every
# subsequent record (100-216) down by one row. The record that used to
# be the last one (row 216) becomes the new last row (217).
$ws.Rows("100:100").Insert()

# Columns that stay constant for every record of this market/product
# need to be (re)written explicitly because Insert() leaves the new row
# completely blank.
$ws.Range("A100").Value = 4
$ws.Range("B100").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C100").Value = "Los Lagos"
$ws.Range("E100").Value = 10
$ws.Range("F100").Value = 100112039
$ws.Range("G100").Value = "Ciboulette"
$ws.Range("H100").Value = "Sin especificar"
$ws.Range("I100").Value = "Primera"
$ws.Range("N100").Value = "$/docena de atados"
$ws.Range("O100").Value = "Región Metropolitana"
$ws.Range("Q100").Value = 3
$ws.Range("R100").Value = "Hortaliza"

# New record's own data.
$ws.Range("D100").Value = 44705
$ws.Range("J100").Value = 240
$ws.Range("K100").Value = 2500
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = 2500
$ws.Range("P100").Value = 833
